# DeepL usage added: insert a new "https://fuckhead.at" frequent-words
# block (5 rows) right after the existing (until now empty) fuckhead.at
# row, followed by a blank separator row, pushing every block that used
# to follow down by 5 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 5 new rows starting at row 3 (row 2 - the original
# "https://fuckhead.at" row with blank word/frequency - stays put and
# simply gets its data filled in below). This pushes the old row 3
# onward ("interstellarrecords" block, "kuprosauwald" block, "röda"
# block, "grgr" block, and their separator rows) down to rows 8-31.
$ws.Rows("3:7").Insert()

# Fill in the word/frequency for the existing fuckhead.at row.
$ws.Range("B2").Value = "Eat"
$ws.Range("C2").Value = 6

# New rows for the rest of the fuckhead.at top-words block.
$ws.Range("A3").Value = "https://fuckhead.at"
$ws.Range("B3").Value = "Immer"
$ws.Range("C3").Value = 6

$ws.Range("A4").Value = "https://fuckhead.at"
$ws.Range("B4").Value = "Kapital"
$ws.Range("C4").Value = 5

$ws.Range("A5").Value = "https://fuckhead.at"
$ws.Range("B5").Value = "Piketty"
$ws.Range("C5").Value = 5

$ws.Range("A6").Value = "https://fuckhead.at"
$ws.Range("B6").Value = "Rich"
$ws.Range("C6").Value = 5

# Row 7 is the blank separator row that follows every word block in
# this sheet. Copying an existing blank separator row (the one that
# used to be row 8, now shifted to row 13, following the
# interstellarrecords.at block) onto row 7 materializes it as a real
# (empty) row instead of leaving a hole in the sheet data.
$ws.Range("A13:C13").Copy($ws.Range("A7:C7"))
